$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4005
$ws.Range("I62").Value = 3980
$ws.Range("J62").Value = 4042.5
$ws.Range("K62").Value = 3980
$ws.Range("L62").Value = 4042.5
$ws.Range("M62").Value = -3356
$ws.Range("N62").Value = -5290.5
$ws.Range("H65").Value = 4005
$ws.Range("I65").Value = 3980
$ws.Range("J65").Value = 4042.5
$ws.Range("K65").Value = 19900
$ws.Range("L65").Value = 20212.5
$ws.Range("M65").Value = -16780
$ws.Range("N65").Value = -26452.5
$ws.Range("H87").Value = 22649
$ws.Range("J87").Value = 22649
$ws.Range("L87").Value = 22649
$ws.Range("N87").Value = -25145
$ws.Range("H90").Value = 22649
$ws.Range("J90").Value = 22649
$ws.Range("L90").Value = 67947
$ws.Range("N90").Value = -80427
$ws.Range("H111").Value = 942.8570999999999
$ws.Range("I111").Value = 899.8
$ws.Range("K111").Value = 2699.4
$ws.Range("M111").Value = 367.6000000000004
$ws.Range("H125").Value = 1234.3158
$ws.Range("J125").Value = 1341.3334
$ws.Range("L125").Value = 12072.0006
$ws.Range("N125").Value = -16992.0006
$ws.Range("H131").Value = 3406.611
$ws.Range("I131").Value = 2772.9167
$ws.Range("J131").Value = 4674
$ws.Range("K131").Value = 8318.750100000001
$ws.Range("L131").Value = 14022
$ws.Range("M131").Value = -3278.750100000001
$ws.Range("N131").Value = -24102
$ws.Range("H135").Value = 894.46155
$ws.Range("I135").Value = 738.9091
$ws.Range("J135").Value = 1750
$ws.Range("K135").Value = 6650.1819
$ws.Range("L135").Value = 15750
$ws.Range("M135").Value = -4115.1819
$ws.Range("N135").Value = -20820

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 19233194
$ws.Range("I2").Value = 31251576
$ws.Range("J2").Value = 3780
$ws.Range("K2").Value = 31251576
$ws.Range("L2").Value = 3780
$ws.Range("M2").Value = -31251463
$ws.Range("N2").Value = -4006
$ws.Range("H32").Value = 6764.4414
$ws.Range("I32").Value = 5748.153
$ws.Range("K32").Value = 5748.153
$ws.Range("M32").Value = -5461.153
$ws.Range("H36").Value = 21848.092
$ws.Range("I36").Value = 3780.25
$ws.Range("K36").Value = 3780.25
$ws.Range("M36").Value = -3434.25
$ws.Range("H110").Value = 2500.3076
$ws.Range("I110").Value = 748.875
$ws.Range("J110").Value = 5302.6
$ws.Range("K110").Value = 748.875
$ws.Range("L110").Value = 5302.6
$ws.Range("M110").Value = 1296.125
$ws.Range("N110").Value = -9392.6
$ws.Range("H116").Value = 19233194
$ws.Range("I116").Value = 31251576
$ws.Range("J116").Value = 3780
$ws.Range("K116").Value = 31251576
$ws.Range("L116").Value = 3780
$ws.Range("M116").Value = -31249282
$ws.Range("N116").Value = -8368

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 19233194
$ws.Range("I3").Value = 31251576
$ws.Range("J3").Value = 3780
$ws.Range("K3").Value = 31251576
$ws.Range("L3").Value = 3780
$ws.Range("M3").Value = -31251462
$ws.Range("N3").Value = -4008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13516893
$ws.Range("I58").Value = 2389.8
$ws.Range("J58").Value = 41672108
$ws.Range("K58").Value = 2389.8
$ws.Range("L58").Value = 41672108
$ws.Range("M58").Value = -2186.8
$ws.Range("N58").Value = -41672514
$ws.Range("H134").Value = 1967.6666
$ws.Range("I134").Value = 866.1177
$ws.Range("J134").Value = 4642.857
$ws.Range("K134").Value = 2598.3531
$ws.Range("L134").Value = 13928.571
$ws.Range("M134").Value = -63.35310000000027
$ws.Range("N134").Value = -18998.571
$ws.Range("H136").Value = 13516893
$ws.Range("I136").Value = 2389.8
$ws.Range("J136").Value = 41672108
$ws.Range("K136").Value = 7169.400000000001
$ws.Range("L136").Value = 125016324
$ws.Range("M136").Value = -4619.400000000001
$ws.Range("N136").Value = -125021424

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 629.2857
$ws.Range("I18").Value = 331
$ws.Range("J18").Value = 1375
$ws.Range("K18").Value = 993
$ws.Range("L18").Value = 4125
$ws.Range("M18").Value = -824
$ws.Range("N18").Value = -4463
$ws.Range("H87").Value = 9996.429
$ws.Range("I87").Value = 1991.6666
$ws.Range("K87").Value = 5974.9998
$ws.Range("M87").Value = -4726.9998
$ws.Range("H90").Value = 9996.429
$ws.Range("I90").Value = 1991.6666
$ws.Range("K90").Value = 17924.9994
$ws.Range("M90").Value = -11684.9994
$ws.Range("H125").Value = 1407.3959
$ws.Range("J125").Value = 1388.9025
$ws.Range("L125").Value = 4166.7075
$ws.Range("N125").Value = -14006.7075
$ws.Range("H130").Value = 2239.4443
$ws.Range("I130").Value = 1788.75
$ws.Range("J130").Value = 2600
$ws.Range("K130").Value = 5366.25
$ws.Range("L130").Value = 7800
$ws.Range("M130").Value = -346.25
$ws.Range("N130").Value = -17840
$ws.Range("H131").Value = 1022.82355
$ws.Range("J131").Value = 1074.5483
$ws.Range("L131").Value = 3223.6449
$ws.Range("N131").Value = -13303.6449

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 18643.715
$ws.Range("J93").Value = 18643.715
$ws.Range("L93").Value = 18643.715
$ws.Range("N93").Value = -22387.715
$ws.Range("H132").Value = 4170.1875
$ws.Range("I132").Value = 4483.375
$ws.Range("J132").Value = 3857
$ws.Range("K132").Value = 13450.125
$ws.Range("L132").Value = 11571
$ws.Range("M132").Value = -10920.125
$ws.Range("N132").Value = -16631

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3847557
$ws.Range("I7").Value = 6250814
$ws.Range("J7").Value = 2345.5
$ws.Range("K7").Value = 6250814
$ws.Range("L7").Value = 2345.5
$ws.Range("M7").Value = -6250702
$ws.Range("N7").Value = -2569.5
$ws.Range("H14").Value = 132123.66
$ws.Range("J14").Value = 7780.5293
$ws.Range("L14").Value = 7780.5293
$ws.Range("N14").Value = -8124.5293
$ws.Range("H61").Value = 2379
$ws.Range("I61").Value = 854.8
$ws.Range("K61").Value = 854.8
$ws.Range("M61").Value = -652.8
$ws.Range("H113").Value = 2379
$ws.Range("I113").Value = 854.8
$ws.Range("K113").Value = 854.8
$ws.Range("M113").Value = 1315.2
$ws.Range("H122").Value = 3398.5334
$ws.Range("I122").Value = 2553.111
$ws.Range("J122").Value = 4666.6665
$ws.Range("K122").Value = 7659.333
$ws.Range("L122").Value = 13999.9995
$ws.Range("M122").Value = -5209.333
$ws.Range("N122").Value = -18899.9995
$ws.Range("H126").Value = 3847557
$ws.Range("I126").Value = 6250814
$ws.Range("J126").Value = 2345.5
$ws.Range("K126").Value = 18752442
$ws.Range("L126").Value = 7036.5
$ws.Range("M126").Value = -18749972
$ws.Range("N126").Value = -11976.5
$ws.Range("H136").Value = 2441388
$ws.Range("I136").Value = 3334523.5
$ws.Range("J136").Value = 5563.636
$ws.Range("K136").Value = 10003570.5
$ws.Range("L136").Value = 16690.908
$ws.Range("M136").Value = -10001020.5
$ws.Range("N136").Value = -21790.908
$ws.Range("H140").Value = 29590.908
$ws.Range("J140").Value = 29590.908
$ws.Range("L140").Value = 29590.908
$ws.Range("N140").Value = -39950.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1491.45
$ws.Range("I113").Value = 398.33334
$ws.Range("K113").Value = 1195.00002
$ws.Range("M113").Value = 974.9999800000001
$ws.Range("H136").Value = 1183.1578
$ws.Range("I136").Value = 550.9259
$ws.Range("J136").Value = 2735
$ws.Range("K136").Value = 1652.7777
$ws.Range("L136").Value = 8205
$ws.Range("M136").Value = 897.2223000000001
$ws.Range("N136").Value = -13305
